# Update the "Users" roster: row 3 team switched from Iowa/Hawkeyes/Big 10
# to Minnesota (keeping the Cyclones/Big 12 values from row 2), and the
# offensive playbook value is corrected to lowercase "flexbone" on both
# data rows. Also clears the stray text-number-format on the Defensive
# Playbook column, widens column D slightly, and moves the selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Minnesota"
$ws.Range("B3").Value = "Cyclones"
$ws.Range("C3").Value = "Big 12"
$ws.Range("F2").Value = "flexbone"
$ws.Range("F3").Value = "flexbone"

$ws.Range("G2:G3").Style = "Normal"

$ws.Columns.Item(4).ColumnWidth = 22.6640625

$ws.Range("A2").Select()
